# edit.ps1 - reproduce the target commit:
#   1) Slide 16's table switches to a different built-in table style
#      ({5B557CB8-9D93-4C94-98E9-D5237B444F2E} -> {B1E98470-FFC1-45BA-A0C4-3D393BE1B9A4}).
#   2) The deck's applied theme ("Integral") is swapped back to the stock
#      "Office Theme" colour palette (the slide master's theme, ppt/theme/theme2.xml).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style swap on slide 16 (the 3rd shape: title, picture, table).
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{B1E98470-FFC1-45BA-A0C4-3D393BE1B9A4}")

# ---------------------------------------------------------------------------
# 2) Re-colour the slide master's theme from "Integral" to the default
#    "Office Theme" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# ---------------------------------------------------------------------------
function Convert-HexToRGBInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Convert-HexToRGBInt $officeColors[$i - 1]
}
